$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item("总计")

function Set-TextValue($range, $text) {
    # Force a numeric-looking string to be stored as text (Excel would
    # otherwise silently coerce e.g. "6.42" to a number on assignment).
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# ---------------------------------------------------------------
# Step 1: duplicate the existing "2022-Q2" sheet to the end of the
# workbook; this copy keeps the original single-fund data untouched
# and becomes the new "2022-Q2" tab.
# ---------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($null, $q2)
$q2copy = $wb.Worksheets.Item($wb.Worksheets.Count)

# ---------------------------------------------------------------
# Step 2: turn the original "2022-Q2" sheet into "2022-Q4" with the
# new quarterly fund-holding data (3 funds). Rename the original
# sheet first so the freshly made copy can claim the "2022-Q2" name.
# ---------------------------------------------------------------
$q4 = $q2
$q4.Name = "2022-Q4"
$q2copy.Name = "2022-Q2"

# headers (row 1) text is unchanged - only restyle to match the
# style used by the "总计" sheet header (s="2").
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

# row 2
$q4.Range("A2").Value = 0
Set-TextValue $q4.Range("B2") "486001"
Set-TextValue $q4.Range("C2") "工银瑞信中国机会全球配置股票（QDII）人民币"
Set-TextValue $q4.Range("D2") "6.42"
Set-TextValue $q4.Range("E2") "93.86"
Set-TextValue $q4.Range("F2") "1.58"
Set-TextValue $q4.Range("G2") "0.1014"
$q4.Range("H2").Value = 8

# row 3
$q4.Range("A3").Value = 1
Set-TextValue $q4.Range("B3") "009562"
Set-TextValue $q4.Range("C3") "工银全球股票（QDII）美元"
Set-TextValue $q4.Range("D3") "6.42"
Set-TextValue $q4.Range("E3") "93.86"
Set-TextValue $q4.Range("F3") "1.58"
Set-TextValue $q4.Range("G3") "0.1014"
$q4.Range("H3").Value = 8

# row 4
$q4.Range("A4").Value = 2
Set-TextValue $q4.Range("B4") "009563"
Set-TextValue $q4.Range("C4") "工银全球股票（QDII）港币"
Set-TextValue $q4.Range("D4") "6.42"
Set-TextValue $q4.Range("E4") "93.86"
Set-TextValue $q4.Range("F4") "1.58"
Set-TextValue $q4.Range("G4") "0.1014"
$q4.Range("H4").Value = 8

# restyle column A (index) to match the "总计" style too (s="2")
# (PasteSpecial copies formats only, so the values set above survive)
$summary.Range("A2").Copy()
$q4.Range("A2:A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# drop the "tabSelected" marker this sheet had as the old active tab
$q2copy.Activate()

# page margins for the new sheet match the workbook-default ones
# used on "总计" (0.75/0.75/1/1/0.5/0.5 in) rather than the old
# 2022-Q2 sheet's (0.7/0.7/0.75/0.75/0.3/0.3 in).
$q4.PageSetup.LeftMargin = 54
$q4.PageSetup.RightMargin = 54
$q4.PageSetup.TopMargin = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------
# Step 3: "总计" sheet - insert a new row 2 with the 2022-Q4 totals,
# pushing the existing 2022-Q2 totals row down to row 3.
# ---------------------------------------------------------------
$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 0.3
$summary.Range("A3").Value = 1
